$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "大有能源"
$ws.Range("B2").Value = "闻泰科技"
$ws.Range("C2").Value = "闻泰科技"
$ws.Range("A3").Value = "神开股份"
$ws.Range("B3").Value = "黄河旋风"
$ws.Range("C3").Value = "大有能源"
$ws.Range("A4").Value = "黄河旋风"
$ws.Range("B4").Value = "工业富联"
$ws.Range("C4").Value = "工业富联"
$ws.Range("A5").Value = "立讯精密"
$ws.Range("B5").Value = "大有能源"
$ws.Range("C5").Value = "合肥城建"
$ws.Range("A6").Value = "三花智控"
$ws.Range("B6").Value = "神开股份"
$ws.Range("C6").Value = "黄河旋风"
$ws.Range("A7").Value = "闻泰科技"
$ws.Range("B7").Value = "立讯精密"
$ws.Range("C7").Value = "神开股份"
$ws.Range("A8").Value = "工业富联"
$ws.Range("B8").Value = "三花智控"
$ws.Range("C8").Value = "安泰科技"
$ws.Range("A9").Value = "中信重工"
$ws.Range("B9").Value = "合肥城建"
$ws.Range("C9").Value = "中信重工"
$ws.Range("A10").Value = "石化机械"
$ws.Range("B10").Value = "天风证券"
$ws.Range("C10").Value = "山河智能"
$ws.Range("A11").Value = "合肥城建"
$ws.Range("B11").Value = "山河智能"
$ws.Range("C11").Value = "华天科技"
$ws.Range("A12").Value = "山河智能"
$ws.Range("B12").Value = "中际旭创"
$ws.Range("C12").Value = "华建集团"
$ws.Range("A13").Value = "天风证券"
$ws.Range("B13").Value = "新易盛"
$ws.Range("C13").Value = "三花智控"
$ws.Range("A14").Value = "中际旭创"
$ws.Range("B14").Value = "山子高科"
$ws.Range("C14").Value = "山东墨龙"
$ws.Range("A15").Value = "德石股份"
$ws.Range("B15").Value = "山东墨龙"
$ws.Range("C15").Value = "山子高科"
$ws.Range("A16").Value = "华天科技"
$ws.Range("B16").Value = "四方达"
$ws.Range("C16").Value = "石化机械"
$ws.Range("A17").Value = "山子高科"
$ws.Range("B17").Value = "中化岩土"
$ws.Range("C17").Value = "立讯精密"
$ws.Range("A18").Value = "山东墨龙"
$ws.Range("B18").Value = "华天科技"
$ws.Range("C18").Value = "上海电力"
$ws.Range("A19").Value = "中化岩土"
$ws.Range("B19").Value = "上海电力"
$ws.Range("C19").Value = "三孚股份"
$ws.Range("A20").Value = "安泰科技"
$ws.Range("B20").Value = "中信重工"
$ws.Range("C20").Value = "福龙马"
$ws.Range("A21").Value = "新易盛"
$ws.Range("B21").Value = "石化机械"
$ws.Range("C21").Value = "楚江新材"
